$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row containing "test121@test.com" / "test123" (row 4) was removed, and a
# brand new row ("testselenium12345@gmail.com" / "test123", with a mailto
# hyperlink on the email cell) was appended at the bottom. Deleting row 4
# shifts every row below it up by one, landing the new entry on row 7 -
# exactly matching the target sheet layout.
$ws.Rows("4:4").Delete() | Out-Null

# Seed A7 with a copy of A2's formatting (the existing hyperlinked cell) so
# the new hyperlink cell lines up with the workbook's existing "link" style
# instead of inventing a new one.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

$ws.Range("A7").Value = "testselenium12345@gmail.com"
$ws.Range("B7").Value = "test123"

# Turn the new email address into a live mailto hyperlink, mirroring A2.
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:testselenium12345@gmail.com", "", "", "testselenium12345@gmail.com") | Out-Null

# Adding the hyperlink resets the cell's style to the engine's generic
# "Hyperlink" style, so reapply the desired formatting afterwards.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

# Rows 4-7 end up with the same (slightly shorter) row height.
$ws.Rows("4:7").RowHeight = 13.8

# Leave the selection where the user's cursor ended up.
$ws.Range("B13").Select() | Out-Null
